$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-15 07:19:00'
$ws.Range('E3').Value = '2026-02-15 07:19:02'
$ws.Range('H3').Value = '92%'
$ws.Range('M3').Value = '-9.1 °C 6:59 TU'
$ws.Range('E4').Value = '2026-02-15 07:19:05'
$ws.Range('J4').Value = '1014.4 hPa'
$ws.Range('L4').Value = '9.4 km/h - 113º 6:44 TU'
$ws.Range('E5').Value = '2026-02-15 07:19:08'
$ws.Range('M5').Value = '-5.6 °C 6:59 TU'
$ws.Range('O5').Value = '-9.4 °C'
$ws.Range('E6').Value = '2026-02-15 07:19:11'
$ws.Range('J6').Value = '1013.8 hPa'
$ws.Range('O6').Value = '4.6 °C'
$ws.Range('E7').Value = '2026-02-15 07:19:13'
$ws.Range('J7').Value = '1013.5 hPa'
$ws.Range('N7').Value = '8.5 °C 6:52 TU'
$ws.Range('E8').Value = '2026-02-15 07:19:16'
$ws.Range('J8').Value = '1013.9 hPa'
$ws.Range('N8').Value = '4.8 °C 6:33 TU'
$ws.Range('E9').Value = '2026-02-15 07:19:19'
$ws.Range('E10').Value = '2026-02-15 07:19:22'
$ws.Range('H10').Value = '75%'
$ws.Range('N10').Value = '-1.2 °C 6:53 TU'
$ws.Range('O10').Value = '3.3 °C'
$ws.Range('E11').Value = '2026-02-15 07:19:24'
$ws.Range('H11').Value = '25%'
$ws.Range('O11').Value = '5.9 °C'
$ws.Range('E12').Value = '2026-02-15 07:19:27'
$ws.Range('E13').Value = '2026-02-15 07:19:30'
$ws.Range('J13').Value = '1015.4 hPa'
$ws.Range('O13').Value = '3.8 °C'
$ws.Range('E14').Value = '2026-02-15 07:19:32'
$ws.Range('O14').Value = '8.9 °C'
$ws.Range('E15').Value = '2026-02-15 07:19:35'
$ws.Range('N15').Value = '7.6 °C 6:46 TU'
$ws.Range('O15').Value = '8.5 °C'
$ws.Range('E16').Value = '2026-02-15 07:19:38'
$ws.Range('H16').Value = '27%'
$ws.Range('E17').Value = '2026-02-15 07:19:41'
$ws.Range('H17').Value = '25%'
$ws.Range('O17').Value = '0.9 °C'
$ws.Range('E18').Value = '2026-02-15 07:19:44'
$ws.Range('H18').Value = '92%'
$ws.Range('J18').Value = '1014.4 hPa'
$ws.Range('K18').Value = '0.0 MJ/m2'
$ws.Range('N18').Value = '-1.1 °C 6:55 TU'
$ws.Range('O18').Value = '1.1 °C'
$ws.Range('E19').Value = '2026-02-15 07:19:46'
$ws.Range('H19').Value = '72%'
$ws.Range('N19').Value = '-1.5 °C 6:51 TU'
$ws.Range('O19').Value = '0.4 °C'
$ws.Range('E20').Value = '2026-02-15 07:19:49'
$ws.Range('E21').Value = '2026-02-15 07:19:52'
$ws.Range('H21').Value = '23%'
$ws.Range('J21').Value = '1014.5 hPa'
$ws.Range('N21').Value = '0.7 °C 6:42 TU'
$ws.Range('O21').Value = '5.1 °C'
$ws.Range('E22').Value = '2026-02-15 07:19:55'
$ws.Range('I22').Value = '0.6 mm'
$ws.Range('N22').Value = '-6.5 °C 6:37 TU'
$ws.Range('E23').Value = '2026-02-15 07:19:57'
$ws.Range('H23').Value = '31%'
$ws.Range('O23').Value = '-6.8 °C'
$ws.Range('E24').Value = '2026-02-15 07:20:00'
$ws.Range('J24').Value = '1017.7 hPa'
$ws.Range('N24').Value = '3.1 °C 6:44 TU'
$ws.Range('O24').Value = '5.1 °C'
$ws.Range('E25').Value = '2026-02-15 07:20:02'
$ws.Range('H25').Value = '51%'
$ws.Range('I25').Value = '14.6 mm'
$ws.Range('O25').Value = '-5.0 °C'
$ws.Range('E26').Value = '2026-02-15 07:20:05'
$ws.Range('E27').Value = '2026-02-15 07:20:07'
$ws.Range('H27').Value = '22%'
$ws.Range('O27').Value = '-2.5 °C'
$ws.Range('E28').Value = '2026-02-15 07:20:10'
$ws.Range('H28').Value = '67%'
$ws.Range('J28').Value = '1014.5 hPa'
$ws.Range('N28').Value = '-1.4 °C 6:51 TU'
$ws.Range('O28').Value = '1.2 °C'
$ws.Range('E29').Value = '2026-02-15 07:20:13'
$ws.Range('H29').Value = '49%'
$ws.Range('O29').Value = '8.0 °C'
$ws.Range('E30').Value = '2026-02-15 07:20:15'
$ws.Range('H30').Value = '51%'
$ws.Range('J30').Value = '1013.3 hPa'
$ws.Range('N30').Value = '4.0 °C 6:59 TU'
$ws.Range('O30').Value = '7.1 °C'
$ws.Range('E31').Value = '2026-02-15 07:20:18'
$ws.Range('J31').Value = '1011.8 hPa'
$ws.Range('E32').Value = '2026-02-15 07:20:21'
$ws.Range('N32').Value = '-0.8 °C 6:42 TU'
$ws.Range('O32').Value = '0.3 °C'
$ws.Range('E33').Value = '2026-02-15 07:20:24'
$ws.Range('J33').Value = '1016.3 hPa'
$ws.Range('N33').Value = '0.1 °C 6:42 TU'
$ws.Range('O33').Value = '2.0 °C'
$ws.Range('E34').Value = '2026-02-15 07:20:26'
$ws.Range('H34').Value = '30%'
$ws.Range('O34').Value = '-1.9 °C'
$ws.Range('E35').Value = '2026-02-15 07:20:29'
$ws.Range('H35').Value = '72%'
$ws.Range('O35').Value = '0.6 °C'
$ws.Range('E36').Value = '2026-02-15 07:20:32'
$ws.Range('J36').Value = '1013.8 hPa'
$ws.Range('N36').Value = '8.3 °C 6:31 TU'
$ws.Range('E37').Value = '2026-02-15 07:20:35'
$ws.Range('H37').Value = '44%'
$ws.Range('J37').Value = '1014.5 hPa'
$ws.Range('N37').Value = '-0.2 °C 6:54 TU'
$ws.Range('O37').Value = '4.2 °C'
$ws.Range('E38').Value = '2026-02-15 07:20:38'
$ws.Range('O38').Value = '2.7 °C'
$ws.Range('E39').Value = '2026-02-15 07:20:40'
$ws.Range('H39').Value = '23%'
$ws.Range('M39').Value = '-4.3 °C 6:50 TU'
$ws.Range('O39').Value = '-5.8 °C'
$ws.Range('E40').Value = '2026-02-15 07:20:43'
$ws.Range('H40').Value = '21%'
$ws.Range('J40').Value = '1014.9 hPa'
$ws.Range('O40').Value = '7.9 °C'
$ws.Range('E41').Value = '2026-02-15 07:20:46'
$ws.Range('J41').Value = '1015.1 hPa'
$ws.Range('E42').Value = '2026-02-15 07:20:49'
$ws.Range('N42').Value = '4.6 °C 6:59 TU'
$ws.Range('O42').Value = '9.1 °C'
$ws.Range('E43').Value = '2026-02-15 07:20:51'
$ws.Range('H43').Value = '73%'
$ws.Range('N43').Value = '-0.7 °C 6:58 TU'
$ws.Range('O43').Value = '2.5 °C'
$ws.Range('E44').Value = '2026-02-15 07:20:54'
$ws.Range('H44').Value = '53%'
$ws.Range('O44').Value = '-7.5 °C'
$ws.Range('E45').Value = '2026-02-15 07:20:57'
$ws.Range('H45').Value = '91%'
$ws.Range('E46').Value = '2026-02-15 07:20:59'
$ws.Range('J46').Value = '1019.3 hPa'
$ws.Range('N46').Value = '7.3 °C 6:57 TU'
